$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 2; this shifts the former rows 2-5
# (新视云, 华为, 满帮, A示例xxx公司) down to rows 3-6.
$ws.Rows.Item(2).Insert()

# Excel's row insert inherits formatting from the row above (the bold/
# bordered header row), which the target file does not want for the new
# data row. Reset it back to the plain/default look used by the other
# data rows before re-applying the one style that IS needed (column A).
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new job entry (创维南京分公司).
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "创维南京分公司"
$ws.Range("C2").Value = "雨花云密城"
$ws.Range("D2").Value = "web后台"
$ws.Range("E2").Value = "Java"
$ws.Range("F2").Value = "09:30"
$ws.Range("G2").Value = "1.5h"
# "995" is purely numeric-looking text in the source file (t="inlineStr"),
# so force it to stay text instead of being auto-coerced to a number --
# same as typing '995 into the cell in the Excel UI.
$ws.Range("H2").Value = "'995"
$ws.Range("I2").Value = "工资八折的10%"
$ws.Range("J2").Value = "1个月工资"
$ws.Range("K2").Value = "不打折"
$ws.Range("L2").Value = "Windows电脑+dell显示器"
$ws.Range("M2").Value = "法定年假"
$ws.Range("N2").Value = "弹性打卡"

# Column A carries the bold / thin-bordered / centered style used by the
# header row and by every other row's "A" (index) cell.
$aCell = $ws.Range("A2")
$aCell.Font.Bold = $true
$aCell.HorizontalAlignment = -4108
$aCell.VerticalAlignment = -4160
$aCell.Borders.LineStyle = 1

# Renumber the sequential index column A for the rows that were shifted down.
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
